$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows 2-4 down to 3-5
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the new occurrence record
$ws.Cells.Item(2, 1).Value = 1851481
$ws.Cells.Item(2, 2).Value = 77595
$ws.Cells.Item(2, 3).Value = "Ovaliderad"
$ws.Cells.Item(2, 4).Value = "LC"
$ws.Cells.Item(2, 5).Value = 6450
$ws.Cells.Item(2, 6).Value = "Skuggblåslav"
$ws.Cells.Item(2, 7).Value = "Hypogymnia vittata"
$ws.Cells.Item(2, 8).Value = "(Ach.) Parrique"
$ws.Cells.Item(2, 16).Value = "Kottorp, Rävbergsmon, Ög"
$ws.Cells.Item(2, 17).Value = 528835.7902055666
$ws.Cells.Item(2, 18).Value = 6489631.821332798
$ws.Cells.Item(2, 19).Value = 100
$ws.Cells.Item(2, 20).Value = "Östergötland"
$ws.Cells.Item(2, 21).Value = "Linköping"
$ws.Cells.Item(2, 22).Value = "Östergötland"
$ws.Cells.Item(2, 23).Value = "Stjärnorp"

# Startdatum / Slutdatum look like dates, so force text so Excel does not
# auto-convert them to date serial numbers.
$dateCell = $ws.Cells.Item(2, 25)
$dateCell.NumberFormat = "@"
$dateCell.Value = "1996-06-05"
$dateCell.Style = "Normal"

$ws.Cells.Item(2, 26).Value = "00:00"

$dateCell2 = $ws.Cells.Item(2, 27)
$dateCell2.NumberFormat = "@"
$dateCell2.Value = "1996-06-05"
$dateCell2.Style = "Normal"

$ws.Cells.Item(2, 28).Value = "00:00"
$ws.Cells.Item(2, 29).Value = "Ek, Tommy, Länsstyrelsens hotartsregister"
$ws.Cells.Item(2, 30).Value = $false
$ws.Cells.Item(2, 31).Value = $false
$ws.Cells.Item(2, 33).Value = $false
$ws.Cells.Item(2, 49).Value = "Linköpings Kommun (hl)"
$ws.Cells.Item(2, 50).Value = "Via Linköpings Kommun (hl)"
$ws.Cells.Item(2, 51).Value = "Linköpings kommuns hotartsdatabas 2012"

# Antal (I) and Bestamningsar (AT) are present but blank in the other rows;
# reproduce that by copying the blank cell from a neighboring row.
$ws.Cells.Item(4, 9).Copy($ws.Cells.Item(2, 9))
$ws.Cells.Item(4, 46).Copy($ws.Cells.Item(2, 46))
